# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G on Sheet1 is labeled "K" (strikeouts). The pitch-by-pitch "Strike#"
# counting previously stored there is replaced with the recalculated
# strikeout totals (K) for each outing, rows 2-63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values for rows 2..63 (column G), in row order.
$kValues = @(
    0, 1, 1, 2, 3, 1, 2, 1, 2, 0,
    1, 0, 1, 2, 0, 1, 1, 1, 1, 0,
    2, 1, 1, 1, 1, 3, 2, 3, 1, 2,
    0, 1, 2, 1, 1, 1, 1, 1, 1, 5,
    2, 2, 2, 0, 5, 0, 2, 2, 0, 1,
    0, 2, 0, 0, 2, 1, 1, 0, 1, 0,
    2, 1
)

$firstRow = 2
for ($i = 0; $i -lt $kValues.Count; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
